$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply "Price" (D) and "Volume(1h)" (E) column updates scraped for this run.
# D-column values are forced to Text format ("@") before assignment so that
# numeric-looking strings (e.g. "560.42", "1.00") keep their exact original
# textual representation (trailing zeros, dot-grouped big numbers, etc.)
# instead of being auto-converted to floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.187.43'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.417.38'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.42'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.00'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.61'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.848.93'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.114.85'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.402.57'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.22'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.52'
$ws.Range('E19').Value = '  +3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '327.47'
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.91'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0772'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.47'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.09'
$ws.Range('E32').Value = '  +7.98%  '
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.43'
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.21'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '324.68'
$ws.Range('E39').Value = '  +3.74%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.56'
$ws.Range('E41').Value = '  -2.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '147.69'
$ws.Range('E42').Value = '  +6.52%  '
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0969'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.80'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.06'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  -0.76%  '
